# Update cryptocurrency price (D) and 1h volume-change (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "82.137.95"
$ws.Range("E2").Value = "  +3.50%  "

$ws.Range("D3").Value = "3.167.83"
$ws.Range("E3").Value = "  -0.27%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.47"
$ws.Range("E5").Value = "  +6.11%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "617.42"
$ws.Range("E6").Value = "  -2.61%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.288"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.15%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.580"
$ws.Range("E9").Value = "  -3.12%  "

$ws.Range("D10").Value = "3.176.41"
$ws.Range("E10").Value = "  -0.10%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.590"
$ws.Range("E11").Value = "  -1.36%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000254"
$ws.Range("E12").Value = "  +0.19%  "

$ws.Range("E13").Value = "  -0.58%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.30"
$ws.Range("E14").Value = "  -1.26%  "

$ws.Range("D15").Value = "3.764.83"
$ws.Range("E15").Value = "  -0.09%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "31.90"
$ws.Range("E16").Value = "  -0.54%  "

$ws.Range("D17").Value = "82.049.13"
$ws.Range("E17").Value = "  +3.44%  "

$ws.Range("D18").Value = "3.180.95"
$ws.Range("E18").Value = "  -0.18%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.90"
$ws.Range("E20").Value = "  -4.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "436.62"
$ws.Range("E21").Value = "  -0.33%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.87"
$ws.Range("E22").Value = "  -5.23%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.11"
$ws.Range("E23").Value = "  -2.26%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.28"
$ws.Range("E24").Value = "  +4.77%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.24"
$ws.Range("E25").Value = "  +9.26%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.89"
$ws.Range("E26").Value = "  +9.97%  "

$ws.Range("D27").Value = "3.333.70"
$ws.Range("E27").Value = "  -0.57%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "76.73"
$ws.Range("E28").Value = "  -0.70%  "

$ws.Range("E29").Value = "  +0.04%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0000120"
$ws.Range("E30").Value = "  -1.74%  "

$ws.Range("E31").Value = "  -0.13%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.97"
$ws.Range("E32").Value = "  -1.31%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "569.84"
$ws.Range("E33").Value = "  +5.50%  "

$ws.Range("E34").Value = "  -2.96%  "

$ws.Range("E35").Value = "  +22.31%  "

$ws.Range("E36").Value = "  -1.60%  "

$ws.Range("E37").Value = "  -2.39%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "22.56"
$ws.Range("E38").Value = "  -1.94%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.19"
$ws.Range("E39").Value = "  +11.96%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.12%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.403"
$ws.Range("E41").Value = "  -1.51%  "

$ws.Range("E42").Value = "  +4.11%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.00"
$ws.Range("E43").Value = "  +11.43%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.98"
$ws.Range("E44").Value = "  +12.67%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "159.21"
$ws.Range("E45").Value = "  -2.96%  "

$ws.Range("E46").Value = "  +0.07%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "185.52"
$ws.Range("E47").Value = "  -3.42%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "44.64"
$ws.Range("E48").Value = "  +2.76%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.31"
$ws.Range("E49").Value = "  -0.52%  "

$ws.Range("E50").Value = "  -3.73%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "25.66"
$ws.Range("E51").Value = "  -0.49%  "
